$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 71.666664
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = -306
$ws.Range("H41").Value = 379.6
$ws.Range("I41").Value = 427.42856
$ws.Range("J41").Value = 337.75
$ws.Range("K41").Value = 427.42856
$ws.Range("L41").Value = 337.75
$ws.Range("M41").Value = 12.57144
$ws.Range("N41").Value = -1217.75
$ws.Range("H99").Value = 908.3125
$ws.Range("I99").Value = 318.57144
$ws.Range("J99").Value = 1367
$ws.Range("K99").Value = 955.71432
$ws.Range("L99").Value = 4101
$ws.Range("M99").Value = 542.28568
$ws.Range("N99").Value = -7097
$ws.Range("H112").Value = 38463108
$ws.Range("J112").Value = 1769
$ws.Range("L112").Value = 5307
$ws.Range("N112").Value = -7523
$ws.Range("H127").Value = 1253.3077
$ws.Range("I127").Value = 527.7143
$ws.Range("J127").Value = 2099.8333
$ws.Range("K127").Value = 1583.1429
$ws.Range("L127").Value = 6299.499899999999
$ws.Range("M127").Value = 3376.8571
$ws.Range("N127").Value = -16219.4999
$ws.Range("H132").Value = 779488.0600000001
$ws.Range("I132").Value = 1651.4314
$ws.Range("J132").Value = 4085293.8
$ws.Range("K132").Value = 4954.2942
$ws.Range("L132").Value = 12255881.4
$ws.Range("M132").Value = -2424.2942
$ws.Range("N132").Value = -12260941.4
$ws.Range("H137").Value = 3126785.5
$ws.Range("I137").Value = 4763149.5
$ws.Range("J137").Value = 2817.7273
$ws.Range("K137").Value = 14289448.5
$ws.Range("L137").Value = 8453.1819
$ws.Range("M137").Value = -14286898.5
$ws.Range("N137").Value = -13553.1819
$ws.Range("H138").Value = 2452667
$ws.Range("I138").Value = 1228.3182
$ws.Range("J138").Value = 6946971
$ws.Range("K138").Value = 3684.9546
$ws.Range("L138").Value = 20840913
$ws.Range("M138").Value = 1455.0454
$ws.Range("N138").Value = -20851193

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 25791.5
$ws.Range("J24").Value = 25791.5
$ws.Range("L24").Value = 25791.5
$ws.Range("N24").Value = -26539.5
$ws.Range("H61").Value = 40080876
$ws.Range("I61").Value = 50050732
$ws.Range("J61").Value = 201450.8
$ws.Range("K61").Value = 50050732
$ws.Range("L61").Value = 201450.8
$ws.Range("M61").Value = -50050520
$ws.Range("N61").Value = -201874.8
$ws.Range("H74").Value = 6001235
$ws.Range("I74").Value = 7844905
$ws.Range("J74").Value = 101490
$ws.Range("K74").Value = 7844905
$ws.Range("L74").Value = 101490
$ws.Range("M74").Value = -7844031
$ws.Range("N74").Value = -103238
$ws.Range("H77").Value = 6001235
$ws.Range("I77").Value = 7844905
$ws.Range("J77").Value = 101490
$ws.Range("K77").Value = 39224525
$ws.Range("L77").Value = 507450
$ws.Range("M77").Value = -39220157
$ws.Range("N77").Value = -516186
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H100").Value = 25791.5
$ws.Range("J100").Value = 25791.5
$ws.Range("L100").Value = 25791.5
$ws.Range("N100").Value = -27955.5
$ws.Range("H122").Value = 4631896.5
$ws.Range("I122").Value = 2399.9473
$ws.Range("J122").Value = 22223982
$ws.Range("K122").Value = 7199.841899999999
$ws.Range("L122").Value = 66671946
$ws.Range("M122").Value = -4749.841899999999
$ws.Range("N122").Value = -66676846
$ws.Range("H124").Value = 30023.2
$ws.Range("J124").Value = 30023.2
$ws.Range("L124").Value = 30023.2
$ws.Range("N124").Value = -39843.2
$ws.Range("H125").Value = 54446.457
$ws.Range("J125").Value = 54446.457
$ws.Range("L125").Value = 54446.457
$ws.Range("N125").Value = -64286.457
$ws.Range("H136").Value = 40080876
$ws.Range("I136").Value = 50050732
$ws.Range("J136").Value = 201450.8
$ws.Range("K136").Value = 150152196
$ws.Range("L136").Value = 604352.3999999999
$ws.Range("M136").Value = -150149646
$ws.Range("N136").Value = -609452.3999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 417.09525
$ws.Range("I94").Value = 279.94116
$ws.Range("K94").Value = 279.94116
$ws.Range("M94").Value = 171.05884
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164
$ws.Range("H134").Value = 1796
$ws.Range("I134").Value = 947.6667
$ws.Range("J134").Value = 2886.7144
$ws.Range("K134").Value = 2843.0001
$ws.Range("L134").Value = 8660.143199999999
$ws.Range("M134").Value = -308.0001000000002
$ws.Range("N134").Value = -13730.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2017.5181
$ws.Range("I31").Value = 951.6491
$ws.Range("J31").Value = 4354.231
$ws.Range("K31").Value = 951.6491
$ws.Range("L31").Value = 4354.231
$ws.Range("M31").Value = -656.6491
$ws.Range("N31").Value = -4944.231
$ws.Range("H32").Value = 1600
$ws.Range("I32").Value = 1600
$ws.Range("K32").Value = 1600
$ws.Range("M32").Value = -1284
$ws.Range("H34").Value = 2017.5181
$ws.Range("I34").Value = 951.6491
$ws.Range("J34").Value = 4354.231
$ws.Range("K34").Value = 951.6491
$ws.Range("L34").Value = 4354.231
$ws.Range("M34").Value = -749.6491
$ws.Range("N34").Value = -4758.231
$ws.Range("H58").Value = 30305260
$ws.Range("I58").Value = 35716200
$ws.Range("J58").Value = 4003
$ws.Range("K58").Value = 35716200
$ws.Range("L58").Value = 4003
$ws.Range("M58").Value = -35715997
$ws.Range("N58").Value = -4409
$ws.Range("H105").Value = 1364.1177
$ws.Range("I105").Value = 1341.3572
$ws.Range("J105").Value = 1470.3334
$ws.Range("K105").Value = 1341.3572
$ws.Range("L105").Value = 1470.3334
$ws.Range("M105").Value = 405.6428000000001
$ws.Range("N105").Value = -4964.3334
$ws.Range("H132").Value = 43350.207
$ws.Range("I132").Value = 29321.805
$ws.Range("J132").Value = 85435.414
$ws.Range("K132").Value = 87965.41500000001
$ws.Range("L132").Value = 256306.242
$ws.Range("M132").Value = -85435.41500000001
$ws.Range("N132").Value = -261366.242
$ws.Range("H134").Value = 20511.965
$ws.Range("I134").Value = 1215.619
$ws.Range("J134").Value = 78401
$ws.Range("K134").Value = 3646.857
$ws.Range("L134").Value = 235203
$ws.Range("M134").Value = -1111.857
$ws.Range("N134").Value = -240273
$ws.Range("H136").Value = 30305260
$ws.Range("I136").Value = 35716200
$ws.Range("J136").Value = 4003
$ws.Range("K136").Value = 107148600
$ws.Range("L136").Value = 12009
$ws.Range("M136").Value = -107146050
$ws.Range("N136").Value = -17109

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 2600
$ws.Range("J42").Value = 2600
$ws.Range("L42").Value = 7800
$ws.Range("N42").Value = -8868
$ws.Range("H131").Value = 860.4
$ws.Range("J131").Value = 1025.6389
$ws.Range("L131").Value = 3076.9167
$ws.Range("N131").Value = -13156.9167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 29750.334
$ws.Range("J92").Value = 29625.5
$ws.Range("L92").Value = 29625.5
$ws.Range("N92").Value = -33369.5
$ws.Range("H97").Value = 2380
$ws.Range("I97").Value = 2380
$ws.Range("K97").Value = 2380
$ws.Range("M97").Value = -1884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 8488.25
$ws.Range("J41").Value = 8488.25
$ws.Range("L41").Value = 8488.25
$ws.Range("N41").Value = -9364.25
$ws.Range("H127").Value = 46333.332
$ws.Range("J127").Value = 46333.332
$ws.Range("L127").Value = 46333.332
$ws.Range("N127").Value = -56253.332
$ws.Range("H132").Value = 27755.45
$ws.Range("I132").Value = 12981.833
$ws.Range("J132").Value = 86849.914
$ws.Range("K132").Value = 38945.499
$ws.Range("L132").Value = 260549.742
$ws.Range("M132").Value = -36415.499
$ws.Range("N132").Value = -265609.742
$ws.Range("H136").Value = 75589.73
$ws.Range("I136").Value = 86161.38
$ws.Range("J136").Value = 67000.25
$ws.Range("K136").Value = 258484.14
$ws.Range("L136").Value = 201000.75
$ws.Range("M136").Value = -255934.14
$ws.Range("N136").Value = -206100.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H113").Value = 2524.625
$ws.Range("J113").Value = 3749.75
$ws.Range("L113").Value = 11249.25
$ws.Range("N113").Value = -15589.25
$ws.Range("H132").Value = 80414.664
$ws.Range("I132").Value = 72318.11
$ws.Range("J132").Value = 101024.09
$ws.Range("K132").Value = 216954.33
$ws.Range("L132").Value = 303072.27
$ws.Range("M132").Value = -214424.33
$ws.Range("N132").Value = -308132.27
$ws.Range("H136").Value = 42223.594
$ws.Range("I136").Value = 31067.758
$ws.Range("J136").Value = 65232.5
$ws.Range("K136").Value = 93203.274
$ws.Range("L136").Value = 195697.5
$ws.Range("M136").Value = -90653.274
$ws.Range("N136").Value = -200797.5

Write-Output "applied 236 cell changes across 8 sheets"